$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cell values in the exact order they were originally authored so the
# shared-strings table is built up in the same sequence as the source edit. ---
$ws.Cells.Item(23, 1).Value = "R22"
$ws.Cells.Item(23, 2).Value = "QA will be handled by Ji while Marina is out from 7/1 to 7/18.   This may slow the process of getting caArray and caIntegrator approved to go to the DEMO/STAGING tier, if we don't get approval before 7/1."
$ws.Cells.Item(24, 1).Value = "R23"
$ws.Cells.Item(25, 1).Value = "R24"
$ws.Cells.Item(24, 2).Value = " We're entering the period when people take vacations.  This may slow communication and impact the pace of work."
$ws.Cells.Item(24, 3).Value = "Mike Hunter"
$ws.Cells.Item(24, 6).Value = "Mike will follow up with team members to find out when they plan to be out and confirm that each team member has identified an appropriate backup to represent them in making decisions or otherwise cover their work while they are away."
$ws.Cells.Item(25, 2).Value = "With caArray 2.5.1 we're going through a new SOP for promotion to DEMO/STAGING tier with new team members (Marina, Mike), so we should expect some initial misses and confusion while we learn the SOP.  This may slow the process of getting approval."

$ws.Cells.Item(23, 4).Value = "Medium"
$ws.Cells.Item(23, 5).Value = "Open"
$ws.Cells.Item(24, 4).Value = "Medium"
$ws.Cells.Item(24, 5).Value = "Open"
$ws.Cells.Item(25, 4).Value = "High"
$ws.Cells.Item(25, 5).Value = "Open"

# --- Formatting: font, borders, wrap text (matches the existing bordered
# "Times New Roman" 12pt risk-table style used through row 22). ---
function Format-RiskCell {
    param($Cell, $Wrap)
    $Cell.Font.Name = "Times New Roman"
    $Cell.Font.Size = 12
    $Cell.Borders.LineStyle = 1
    $Cell.WrapText = $Wrap
}

Format-RiskCell $ws.Cells.Item(23, 1) $false
Format-RiskCell $ws.Cells.Item(23, 2) $true
Format-RiskCell $ws.Cells.Item(23, 3) $true
Format-RiskCell $ws.Cells.Item(23, 4) $true
Format-RiskCell $ws.Cells.Item(23, 5) $false
Format-RiskCell $ws.Cells.Item(23, 6) $false

Format-RiskCell $ws.Cells.Item(24, 1) $false
Format-RiskCell $ws.Cells.Item(24, 2) $true
Format-RiskCell $ws.Cells.Item(24, 3) $true
Format-RiskCell $ws.Cells.Item(24, 4) $true
Format-RiskCell $ws.Cells.Item(24, 5) $false
Format-RiskCell $ws.Cells.Item(24, 6) $true

Format-RiskCell $ws.Cells.Item(25, 1) $false
Format-RiskCell $ws.Cells.Item(25, 2) $true
Format-RiskCell $ws.Cells.Item(25, 3) $true
Format-RiskCell $ws.Cells.Item(25, 4) $true
Format-RiskCell $ws.Cells.Item(25, 5) $false
Format-RiskCell $ws.Cells.Item(25, 6) $false

# --- Row heights (wrapped text row heights from the source edit). ---
$ws.Rows(23).RowHeight = 90
$ws.Rows(24).RowHeight = 60
$ws.Rows(25).RowHeight = 105

# --- Selection / view state matches the post-edit workbook. ---
$ws.Range("F24").Select()
